$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 93 ---
$ws.Range("E93").Value = 2
$ws.Range("H93").Value = 44371
$ws.Range("H92").Copy()
$ws.Range("H93").PasteSpecial(-4122)
$ws.Range("I93").Value = 45863

# --- New row 100 ---
$ws.Range("A100").Value = 713
$ws.Range("B100").Value = "Subarray Product Less Than K"
$ws.Range("C100").Value = "#array #two-pointers #sliding-window #核心 #count-subarrays"
$ws.Range("D100").Value = "medium"
$ws.Range("E100").Value = 0
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 50
$ws.Range("H100").Value = 45864
$ws.Range("I93").Copy()
$ws.Range("H100").PasteSpecial(-4122)
$ws.Range("I100").Value = 45864
$ws.Range("I93").Copy()
$ws.Range("I100").PasteSpecial(-4122)
$ws.Range("J100").Value = "template: slide window"
$ws.Rows("100").RowHeight = 68

# --- New row 101 ---
$ws.Range("A101").Value = 395
$ws.Range("B101").Value = "Longest Substring with At Least K Repeating Characters"
$ws.Range("C101").Value = "#hash-table #string #sliding-window #divide-and-conquer  #重点 "
$ws.Range("D101").Value = "medium"
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 1
$ws.Range("G101").Value = 40
$ws.Range("H101").Value = 45864
$ws.Range("I93").Copy()
$ws.Range("H101").PasteSpecial(-4122)
$ws.Range("I101").Value = 45864
$ws.Range("I93").Copy()
$ws.Range("I101").PasteSpecial(-4122)
$ws.Range("J101").Value = "template: slide window"
$ws.Rows("101").RowHeight = 84

# --- New row 102 ---
$ws.Range("A102").Value = 3480
$ws.Range("B102").Value = "Maximize Subarrays After Removing One Conflicting Pair"
$ws.Range("C102").Value = "#array #segment-tree #enumeration #prefix-sum "
$ws.Range("D102").Value = "hard"
$ws.Range("E102").Value = 0
$ws.Range("F102").Value = 1
$ws.Range("G102").Value = 50
$ws.Range("H102").Value = 45865
$ws.Range("I93").Copy()
$ws.Range("H102").PasteSpecial(-4122)
$ws.Range("I102").Value = 45865
$ws.Range("I93").Copy()
$ws.Range("I102").PasteSpecial(-4122)
$ws.Rows("102").RowHeight = 51

# --- New row 103 ---
$ws.Range("A103").Value = 2210
$ws.Range("B103").Value = "Count Hills and Valleys in an Array"
$ws.Range("C103").Value = "#array"
$ws.Range("D103").Value = "ezsy"
$ws.Range("E103").Value = 1
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 10
$ws.Range("H103").Value = 45865
$ws.Range("I93").Copy()
$ws.Range("H103").PasteSpecial(-4122)
$ws.Range("I103").Value = 45865
$ws.Range("I93").Copy()
$ws.Range("I103").PasteSpecial(-4122)
$ws.Rows("103").RowHeight = 34

# --- New row 104 ---
$ws.Range("A104").Value = 2044
$ws.Range("B104").Value = "Count Number of Maximum Bitwise-OR Subsets"
$ws.Range("C104").Value = "#memoization #backtracking #recursive #dynamic-programming "
$ws.Range("D104").Value = "medium"
$ws.Range("E104").Value = 0
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 45
$ws.Range("H104").Value = 45866
$ws.Range("I93").Copy()
$ws.Range("H104").PasteSpecial(-4122)
$ws.Range("I104").Value = 45866
$ws.Range("I93").Copy()
$ws.Range("I104").PasteSpecial(-4122)
$ws.Rows("104").RowHeight = 84

# --- Selection / view ---
$null = $ws.Range("B107").Select()
